$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 94.22446084399637
$ws.Range("H2").Value = 98.82114311930236
$ws.Range("I2").Value = 90.54662346319651

$ws.Range("G3").Value = 99.40169005610092
$ws.Range("H3").Value = 97.23717927663057
$ws.Range("I3").Value = 97.44278969460754

$ws.Range("G4").Value = 98.20421142833979
$ws.Range("H4").Value = 98.70888207145933
$ws.Range("I4").Value = 96.72664470863005

$ws.Range("G5").Value = 97.84177238794999
$ws.Range("H5").Value = 98.22390845754153
$ws.Range("I5").Value = 94.70304416476111

$ws.Range("G6").Value = 98.52366317681302
$ws.Range("H6").Value = 98.55556892001933
$ws.Range("I6").Value = 94.73406432600048

$wb.Save()
